$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.680.59'
$ws.Range('E2').Value = '  -2.22%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.016.84'
$ws.Range('E3').Value = '  -4.33%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').Value = '  +0.58%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '332.83'
$ws.Range('E5').Value = '  -3.64%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.011'
$ws.Range('E6').Value = '  +0.42%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5064'
$ws.Range('E7').Value = '  -3.06%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4266'
$ws.Range('E8').Value = '  -3.93%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '53.95'
$ws.Range('E9').Value = '  -0.82%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.09249'
$ws.Range('E10').Value = '  -2.37%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.130'
$ws.Range('E11').Value = '  -3.73%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.64'
$ws.Range('E12').Value = '  -5.79%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.145'
$ws.Range('E13').Value = '  -6.67%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.575'
$ws.Range('E14').Value = '  -4.88%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.946.35'
$ws.Range('E15').Value = '  -8.08%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '95.71'
$ws.Range('E16').Value = '  -5.89%  '

# Row 17
$ws.Range('E17').Value = '  +0.51%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001128'
$ws.Range('E18').Value = '  -3.08%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06677'
$ws.Range('E19').Value = '  -0.85%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.00'
$ws.Range('E20').Value = '  -6.10%  '

# Row 21
$ws.Range('E21').Value = '  +0.39%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.026'
$ws.Range('E22').Value = '  -4.67%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '29.707.95'
$ws.Range('E23').Value = '  -2.26%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.06'
$ws.Range('E24').Value = '  -4.65%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.270'
$ws.Range('E25').Value = '  -2.11%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.90'
$ws.Range('E26').Value = '  -2.09%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.87'
$ws.Range('E27').Value = '  -5.27%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.415'
$ws.Range('E28').Value = '  -7.18%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.359'
$ws.Range('E29').Value = '  -7.23%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '128.92'
$ws.Range('E30').Value = '  -3.37%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.064'
$ws.Range('E31').Value = '  -7.26%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.602'
$ws.Range('E32').Value = '  -8.99%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.1001'
$ws.Range('E33').Value = '  -5.14%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.895'
$ws.Range('E34').Value = '  -5.85%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.805'
$ws.Range('E35').Value = '  -3.23%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.636'
$ws.Range('E36').Value = '  -8.50%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02477'
$ws.Range('E37').Value = '  -6.23%  '

# Row 38
$ws.Range('E38').Value = '  -1.56%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06418'
$ws.Range('E39').Value = '  -5.66%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6606'
$ws.Range('E40').Value = '  -6.13%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.86'
$ws.Range('E41').Value = '  -5.53%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.2092'
$ws.Range('E42').Value = '  -6.09%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.010'
$ws.Range('E43').Value = '  +0.41%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6394'
$ws.Range('E44').Value = '  -6.54%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.65'
$ws.Range('E45').Value = '  -5.61%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.225'
$ws.Range('E46').Value = '  -5.90%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.293'
$ws.Range('E47').Value = '  -4.57%  '

# Row 48
$ws.Range('E48').Value = '  -2.96%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.07032'
$ws.Range('E49').Value = '  -2.77%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00000000325'
$ws.Range('E50').Value = '  -6.14%  '

# Row 51
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.143'
$ws.Range('E51').Value = '  -4.99%  '

